# Update the "All" label to "Combined" in column A for the summary rows
# of the race/ethnicity table (rows 2, 5, 8, 11, 14, 17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 5, 8, 11, 14, 17)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = "Combined"
}
